$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = $tr.Text.Replace("First Prioritization Grid.", "First Prioritization Grid")
